{"js": "// Remove the paragraph that contains the hyperlink to\n// https://youtu.be/cCOz5xle820 (the \"trailer moved\" link), located right\n// after the \"Link video sustentaci\u00f3n\" paragraph and right before the\n// trailing bookmark paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = \"https://youtu.be/cCOz5xle820\";\n\nlet paragraphToDelete = null;\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text && paragraph.text.indexOf(target) !== -1) {\n    paragraphToDelete = paragraph;\n    break;\n  }\n}\n\nif (paragraphToDelete) {\n  paragraphToDelete.delete();\n  await context.sync();\n}\n", "ps1": "# The trailer video was moved, so the paragraph holding the hyperlink\n# https://youtu.be/cCOz5xle820 (right after the \"Link video sustentaci\u00f3n\"\n# paragraph and right before the trailing bookmark paragraph) is removed\n# entirely, paragraph mark included.\n\n$d = $word.ActiveDocument\n$target = \"https://youtu.be/cCOz5xle820\"\n\n$found = $false\n\n# Preferred idiom: locate the link text with Find, then drop the whole\n# paragraph that contains it (this also removes its paragraph mark, merging\n# it away just like Word does when you select the paragraph and hit Delete).\n$rng = $d.Content\nif ($rng.Find.Execute($target)) {\n    $para = $rng.Paragraphs(1)\n    $para.Range.Delete()\n    $found = $true\n}\n\nif (-not $found) {\n    # Fallback: walk the paragraphs and match on text.\n    foreach ($p in @($d.Paragraphs)) {\n        $r = $p.Range\n        if ($r.Text -like \"*$target*\") {\n            $r.Delete()\n            break\n        }\n    }\n}\n"}
